# B6-PowerPoint.pptx — 28 Jun 2020 commit
#
# 1) Three tables (on the slides holding the "Profitability / Liquidity /
#    gearing" ratio tables) get their table style switched from the
#    Google-Slides-exported custom style {8AFFED38-...} to the built-in
#    "No Style, Table Grid" style {0655719B-...}.
# 2) The deck's theme (ppt/theme/theme1.xml, "Integral" / "Red Violet")
#    is recoloured to the stock Office theme palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three ratio tables -------------------------------
$tableSlides = 14, 15, 16
foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{0655719B-AD0B-44C6-B916-60D8798AD9D4}")
    }
}

# --- 2) Recolour the theme to the default Office palette --------------
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# msoThemeColorDark1=1 .. msoThemeColorFollowedHyperlink=12
$officeColors = @{
    1  = (RGBVal 0x00 0x00 0x00)   # Dark 1
    2  = (RGBVal 0xFF 0xFF 0xFF)   # Light 1
    3  = (RGBVal 0x44 0x54 0x6A)   # Dark 2
    4  = (RGBVal 0xE7 0xE6 0xE6)   # Light 2
    5  = (RGBVal 0x5B 0x9B 0xD5)   # Accent 1
    6  = (RGBVal 0xED 0x7D 0x31)   # Accent 2
    7  = (RGBVal 0xA5 0xA5 0xA5)   # Accent 3
    8  = (RGBVal 0xFF 0xC0 0x00)   # Accent 4
    9  = (RGBVal 0x44 0x72 0xC4)   # Accent 5
    10 = (RGBVal 0x70 0xAD 0x47)   # Accent 6
    11 = (RGBVal 0x05 0x63 0xC1)   # Hyperlink
    12 = (RGBVal 0x95 0x4F 0x72)   # Followed Hyperlink
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i]
}
